$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.956.15'
$ws.Range("E2").Value = '  -4.37%  '
$ws.Range("D3").Value = '2.233.29'
$ws.Range("E3").Value = '  -4.94%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.98'
$ws.Range("E5").Value = '  +1.60%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.628'
$ws.Range("E6").Value = '  -5.99%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '68.31'
$ws.Range("E7").Value = '  -6.95%  '
$ws.Range("E8").Value = '  +0.11%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.554'
$ws.Range("E9").Value = '  -7.48%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0978'
$ws.Range("E10").Value = '  -3.38%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '58.64'
$ws.Range("E11").Value = '  -2.12%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '35.80'
$ws.Range("E12").Value = '  +8.62%  '
$ws.Range("E13").Value = '  -2.83%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.71'
$ws.Range("E14").Value = '  -7.82%  '
$ws.Range("D15").Value = '2.570.48'
$ws.Range("E15").Value = '  -4.75%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.87'
$ws.Range("E16").Value = '  -7.98%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.857'
$ws.Range("E17").Value = '  -5.54%  '
$ws.Range("D18").Value = '2.231.17'
$ws.Range("E18").Value = '  -5.18%  '
$ws.Range("D19").Value = '42.000.40'
$ws.Range("E19").Value = '  -4.16%  '
$ws.Range("D20").Value = '0.0₃0961'
$ws.Range("E20").Value = '  -6.17%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.90'
$ws.Range("E21").Value = '  -7.00%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.19'
$ws.Range("E22").Value = '  -7.45%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '235.11'
$ws.Range("E23").Value = '  -6.63%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.07'
$ws.Range("E24").Value = '  +11.80%  '
$ws.Range("E25").Value = '  -0.08%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.63'
$ws.Range("E26").Value = '  -3.72%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.45'
$ws.Range("E27").Value = '  -1.74%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.25'
$ws.Range("E28").Value = '  -2.39%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.95'
$ws.Range("E29").Value = '  -4.87%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '171.75'
$ws.Range("E30").Value = '  -2.71%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '20.46'
$ws.Range("E31").Value = '  -8.06%  '
$ws.Range("E32").Value = '  -4.80%  '
$ws.Range("E33").Value = '  -5.78%  '
$ws.Range("E34").Value = '  -4.38%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.22'
$ws.Range("E35").Value = '  -2.71%  '
$ws.Range("E36").Value = '  -7.90%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.80'
$ws.Range("E37").Value = '  +0.02%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '23.67'
$ws.Range("E38").Value = '  +25.68%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0282'
$ws.Range("E39").Value = '  +3.84%  '
$ws.Range("E40").Value = '  -3.36%  '
$ws.Range("B41").Value = 'MultiversX'
$ws.Range("C41").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '66.79'
$ws.Range("E41").Value = '  +1.68%  '
$ws.Range("B42").Value = 'THORChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.84'
$ws.Range("E42").Value = '  -8.78%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '9.05'
$ws.Range("E43").Value = '  -1.91%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.92'
$ws.Range("E44").Value = '  -12.60%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.101'
$ws.Range("E45").Value = '  -3.86%  '
$ws.Range("E46").Value = '  -3.09%  '
$ws.Range("E47").Value = '  +0.30%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.48'
$ws.Range("E48").Value = '  +5.03%  '
$ws.Range("E49").Value = '  -2.74%  '
$ws.Range("B50").Value = 'HuobiToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.82'
$ws.Range("E50").Value = '  -2.21%  '
$ws.Range("E51").Value = '  -4.29%  '
